$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E4 value from 1900 to 1800
$ws.Range("E4").Value = 1800

# Add new row 17
$ws.Range("A17").Value = "NONEXIST-SKU-9999-XXX"
$ws.Range("B17").Value = "Ghost"
$ws.Range("C17").Value = "Phantom"
$ws.Range("D17").Value = 99999
$ws.Range("E17").Value = 2023
$ws.Range("F17").Value = "black"
